$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Window view size (bookViews / workbookView) ---
$win = $excel.ActiveWindow
$win.Width = 20700
$win.Height = 9450

# --- Copyright year text update (B3) ---
$ws.Range("B3").Value = "Copyright @2015 - 2023"

# --- Quarter headers row 8 (D8:H8): drop oldest quarter, shift left, append new quarter ---
$ws.Range("D8").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("E8").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("F8").Value = "فصل اول منتهی به 1401/03"
$ws.Range("G8").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("H8").Value = "فصل سوم منتهی به 1401/09"

# --- Publish date headers row 9 (D9:H9) ---
# NB: plain "YYYY-MM-DD" text (no trailing suffix) is auto-parsed into a
# date serial by the engine's smart-entry heuristic, same as real Excel
# would on literal keystrokes into a General-formatted cell. A leading
# apostrophe forces these two through as literal text, matching the
# shared-string type the source file expects.
$ws.Range("D9").Value = "1401-11-02 (2)"
$ws.Range("E9").Value = "1401-11-02 (7)"
$ws.Range("F9").Value = "'1401-04-29"
$ws.Range("G9").Value = "1401-09-16 (3)"
$ws.Range("H9").Value = "'1401-11-02"

# --- Shift financial data in rows 11-27 (excluding blank rows 15/23) one column to the left,
#     dropping the old quarter's data (old column D) and filling in new data for column H ---
$dataRows = @(11,12,13,14,16,17,18,19,20,21,22,24,25,26,27)

$newValues = @{
    11 = @(5415036, 7328558, 5730428, 9439640, 10285175)
    12 = @(-4248442, -5649902, -4286727, -6731654, -6644712)
    13 = @(1166594, 1678656, 1443701, 2707986, 3640463)
    14 = @(-474931, -845696, -714296, -635801, -1260601)
    16 = @(13956, 28773, -4108, -50542, -8354)
    17 = @(705619, 861733, 725297, 2021643, 2371508)
    18 = @(-380317, -461026, -559986, -515326, -821021)
    19 = @(-494501, 648132, 49622, 688653, -646718)
    20 = @(-169199, 1048839, 214933, 2194970, 903769)
    21 = @(7795, -37203, -10332, -19317, -99281)
    22 = @(-161404, 1011636, 204601, 2175653, 804488)
    24 = @(-161404, 1011636, 204601, 2175653, 804488)
    25 = @(-81, 506, 102, 1088, 402)
    26 = @(2000000, 2000000, 2000000, 2000000, 2000000)
    27 = @(-81, 506, 102, 1088, 402)
}

$cols = @("D", "E", "F", "G", "H")

foreach ($r in $dataRows) {
    $vals = $newValues[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}

# --- Column widths: D/E grouped at 29 + F at 31 + G/H at 29 previously;
#     now D alone at 29, E alone at 31, F/G/H grouped at 29 ---
$ws.Range("D1:D28").ColumnWidth = 29
$ws.Range("E1:E28").ColumnWidth = 31
$ws.Range("F1:H28").ColumnWidth = 29
